$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.99999999106689519
$ws.Range("A2").Value = 0.99460444540042214
$ws.Range("A3").Value = 0.97448129466878908
$ws.Range("A4").Value = 0.96549876239454513
$ws.Range("A5").Value = 0.95695159974188937
$ws.Range("A6").Value = 0.93609948796318432
$ws.Range("A7").Value = 0.9345123716502417
$ws.Range("A8").Value = 0.93308811809745784
$ws.Range("A9").Value = 0.9345430067893854
$ws.Range("A10").Value = 0.93697339299768467
$ws.Range("A11").Value = 0.93747930686102388
$ws.Range("A12").Value = 0.93869973838003529
$ws.Range("A13").Value = 0.94656466313588794
$ws.Range("A14").Value = 0.95094756370256595
$ws.Range("A15").Value = 0.9563801005110677
$ws.Range("A16").Value = 0.95387396500027299
$ws.Range("A17").Value = 0.95016640993388135
$ws.Range("A18").Value = 0.94905752866523152
$ws.Range("A19").Value = 0.99143136485342021
$ws.Range("A20").Value = 0.98431431973501338
$ws.Range("A21").Value = 0.98291584006905897
$ws.Range("A22").Value = 0.98165133429005347
$ws.Range("A23").Value = 0.96744582937803181
$ws.Range("A24").Value = 0.9544242567587442
$ws.Range("A25").Value = 0.94796715130385212
$ws.Range("A26").Value = 0.93301677461528798
$ws.Range("A27").Value = 0.9289682050254302
$ws.Range("A28").Value = 0.91119060964393928
$ws.Range("A29").Value = 0.89865779806652801
$ws.Range("A30").Value = 0.89345030356540267
$ws.Range("A31").Value = 0.88926804583699282
$ws.Range("A32").Value = 0.8875887647215519
$ws.Range("A33").Value = 0.88706876776429056
